# Update "F" column ("想去人数" / want-to-go count) values on the
# "展览", "演出" and "全部类型" sheets, per the upstream scraper refresh.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @{
    "F3"  = 78
    "F4"  = 9552
    "F5"  = 629
    "F6"  = 164
    "F7"  = 299
    "F8"  = 360
    "F9"  = 417
    "F11" = 196
    "F13" = 451
    "F14" = 12183
    "F21" = 176
    "F24" = 2731
    "F25" = 2101
    "F28" = 56
    "F29" = 2150
    "F30" = 1014
    "F31" = 4207
    "F32" = 3665
    "F33" = 587
    "F34" = 2627
    "F36" = 31
    "F42" = 429
    "F43" = 543
    "F46" = 231
    "F49" = 144
}

$sheet2Updates = @{
    "F8"  = 56
    "F9"  = 44
    "F17" = 27
    "F20" = 187
    "F22" = 9
    "F24" = 76
}

$sheet4Updates = @{
    "F3"  = 56
    "F5"  = 78
    "F6"  = 9552
    "F7"  = 629
    "F8"  = 44
    "F9"  = 164
    "F10" = 299
    "F11" = 360
    "F12" = 417
    "F14" = 196
    "F15" = 451
    "F16" = 12183
    "F22" = 176
    "F25" = 2731
    "F26" = 2101
    "F29" = 56
    "F30" = 2150
    "F31" = 1014
    "F32" = 4207
    "F33" = 3665
    "F34" = 587
    "F35" = 2627
    "F41" = 429
    "F43" = 543
    "F46" = 231
    "F49" = 144
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

$ws2 = $wb.Worksheets.Item("演出")
foreach ($addr in $sheet2Updates.Keys) {
    $ws2.Range($addr).Value = $sheet2Updates[$addr]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}
